$d = $word.ActiveDocument

# 1. Fix table reference: "Table A3" -> "Table A2"
$d.Content.Find.Execute("Table A3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Table A2", 2)
